$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'40.918.48"
$ws.Range("E2").Value = "  -2.26%  "
$ws.Range("D3").Value = "'2.174.92"
$ws.Range("E3").Value = "  -2.91%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'250.01"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").Value = "'0.620"
$ws.Range("E6").Value = "  -2.18%  "
$ws.Range("D7").Value = "'66.85"
$ws.Range("E7").Value = "  -7.48%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.568"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "'58.66"
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("D11").Value = "'0.0930"
$ws.Range("E11").Value = "  -4.37%  "
$ws.Range("D12").Value = "'35.30"
$ws.Range("E12").Value = "  -16.47%  "
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("D14").Value = "'6.93"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("D15").Value = "'2.494.16"
$ws.Range("E15").Value = "  -3.04%  "
$ws.Range("D16").Value = "'0.859"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "'14.15"
$ws.Range("E17").Value = "  -6.45%  "
$ws.Range("D18").Value = "'2.175.58"
$ws.Range("E18").Value = "  -2.83%  "
$ws.Range("D19").Value = "'40.927.12"
$ws.Range("E19").Value = "  -2.01%  "
$ws.Range("D20").Value = "'0.0₃0942"
$ws.Range("E20").Value = "  -2.73%  "
$ws.Range("D21").Value = "'6.09"
$ws.Range("E21").Value = "  -1.99%  "
$ws.Range("D22").Value = "'71.55"
$ws.Range("E22").Value = "  -2.57%  "
$ws.Range("D23").Value = "'230.36"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("D24").Value = "'2.07"
$ws.Range("E24").Value = "  -8.95%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'11.34"
$ws.Range("E26").Value = "  +11.89%  "
$ws.Range("D27").Value = "'3.71"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").Value = "'2.42"
$ws.Range("E28").Value = "  -2.81%  "
$ws.Range("D29").Value = "'2.14"
$ws.Range("E29").Value = "  -2.85%  "
$ws.Range("D30").Value = "'167.83"
$ws.Range("E30").Value = "  -2.36%  "
$ws.Range("D31").Value = "'20.26"
$ws.Range("E31").Value = "  -2.66%  "
$ws.Range("D32").Value = "'0.122"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").Value = "'5.61"
$ws.Range("E33").Value = "  +2.91%  "
$ws.Range("D34").Value = "'0.0749"
$ws.Range("E34").Value = "  +3.47%  "
$ws.Range("D35").Value = "'0.122"
$ws.Range("E35").Value = "  -2.73%  "
$ws.Range("D36").Value = "'4.09"
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("D37").Value = "'4.55"
$ws.Range("E37").Value = "  -3.59%  "
$ws.Range("D38").Value = "'25.27"
$ws.Range("E38").Value = "  -5.54%  "
$ws.Range("D39").Value = "'0.0300"
$ws.Range("E39").Value = "  +7.58%  "
$ws.Range("D40").Value = "'2.18"
$ws.Range("E40").Value = "  -4.57%  "
$ws.Range("B41").Value = "THORChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D41").Value = "'5.51"
$ws.Range("E41").Value = "  -8.62%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "'11.67"
$ws.Range("E42").Value = "  +0.56%  "
$ws.Range("D43").Value = "'61.84"
$ws.Range("E43").Value = "  -9.59%  "
$ws.Range("D44").Value = "'4.80"
$ws.Range("E44").Value = "  -4.99%  "
$ws.Range("D45").Value = "'0.193"
$ws.Range("E45").Value = "  -10.71%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'8.59"
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("B47").Value = "BinanceUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'1.18"
$ws.Range("E48").Value = "  +5.13%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0986"
$ws.Range("E49").Value = "  -3.22%  "
$ws.Range("D50").Value = "'1.15"
$ws.Range("E50").Value = "  -3.36%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").Value = "'2.71"
$ws.Range("E51").Value = "  -0.41%  "
